$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity for RESISTOR 4700 (row 28) from 1 to 2
$ws.Range("B28").Value = 2

# Scroll the view so that A23 is the top-left visible cell
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
